$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ACS")

# The sheet already has 5 blank placeholder rows (46-50) before the
# "driving_alone / commute" block. The new "digital divide" variables need
# 6 rows (46-51), so insert 2 more blank rows ahead of the old row 50 —
# this pushes the old rows 50-57 down to 52-59, matching the target layout.
$ws.Range("A50:A51").EntireRow.Insert()

# Rows 46-49 already carry the right cell formatting (blank placeholder
# cells styled like the rest of the table) - just fill in the values, in
# the same left-to-right order the workbook's author used.
$ws.Range("B46").Value2 = "CINETHH"
$ws.Range("C46").Value2 = "digital divide"
$ws.Range("E46").Value2 = "access to internet "

$ws.Range("B47").Value2 = "CILAPTOP"
$ws.Range("C47").Value2 = "digital divide"
$ws.Range("E47").Value2 = "laptop, desktop, or notebook computer"

$ws.Range("B48").Value2 = "CISMRTPHN"
$ws.Range("C48").Value2 = "digital divide"
$ws.Range("E48").Value2 = "smartphone"

$ws.Range("B49").Value2 = "CITABLET"
$ws.Range("C49").Value2 = "digital divide"
$ws.Range("E49").Value2 = "tablet or other portable wireless computer"
$ws.Range("D49").Value2 = "unsure if this overlaps with laptop"

# Rows 50-51 are the freshly-inserted rows: they inherited row 49's format
# on insert, so restyle them to match the rest of the "blank separator
# row above" style block (same look as rows 33-36 / the row that follows).
foreach ($r in 50,51) {
    $ws.Range("A$r").HorizontalAlignment = -4108
    $ws.Range("A$r").VerticalAlignment = -4107
    $ws.Range("A$r").WrapText = $true

    $ws.Range("B$r").HorizontalAlignment = -4108
    $ws.Range("B$r").VerticalAlignment = -4107
    $ws.Range("B$r").WrapText = $false

    $ws.Range("C$r").HorizontalAlignment = -4108
    $ws.Range("C$r").VerticalAlignment = -4107
    $ws.Range("C$r").WrapText = $true

    $ws.Range("D$r").HorizontalAlignment = 1
    $ws.Range("D$r").VerticalAlignment = -4107
    $ws.Range("D$r").WrapText = $true
}

$ws.Range("B50").Value2 = "CIHISPEED"
$ws.Range("C50").Value2 = "digital divide"
$ws.Range("E50").Value2 = "broadband (high speed) internet service such as cable, fiber optic, or DSL service"

$ws.Range("B51").Value2 = "CIDIAL"
$ws.Range("C51").Value2 = "digital divide"
$ws.Range("E51").Value2 = "dial-up service"
